$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update municipality names with correct accents / formatting.
# Only the displayed text changes; row positions and all numeric data stay the same.
$ws.Range("A2").Value  = "C.A.B.A."
$ws.Range("A6").Value  = "SAN MARTÍN"
$ws.Range("A7").Value  = "VICENTE LÓPEZ"
$ws.Range("A10").Value = "ITUZAINGÓ"
$ws.Range("A11").Value = "MORÓN"
$ws.Range("A14").Value = "LANÚS"
$ws.Range("A16").Value = "ALMIRANTE BROWN"
$ws.Range("A17").Value = "ESTEBAN ECHEVERRÍA"
$ws.Range("A21").Value = "JOSÉ C. PAZ"
$ws.Range("A28").Value = "PRESIDENTE PERÓN"

# Move the active selection to A29 (below the last data row), matching the
# saved sheet view state in the workbook.
$ws.Range("A29").Select()

# Record the workbook window position, matching the saved bookViews state
# (xWindow/yWindow in workbook.xml).
$excel.ActiveWindow.Left = 9640
$excel.ActiveWindow.Top = 40
